$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns holding date-like text values ("YYYY-MM-DD") that must stay plain
# text (otherwise Excel auto-converts the assigned string into a real date
# serial number). Briefly force a Text number format while writing, then
# drop the format again so the cell's style matches the untouched original.
$textCols = @(25, 27)  # Y, AA

$lastCol = 51  # column AY (last used column on the sheet)

for ($c = 1; $c -le $lastCol; $c++) {
    $cell2 = $ws.Cells.Item(2, $c)
    $cell3 = $ws.Cells.Item(3, $c)

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    # Nothing to do if both sides already hold the same (blank) content -
    # writing "" back would delete an otherwise-present empty text cell.
    if ($v2 -eq $v3) {
        continue
    }

    $isTextCol = $textCols -contains $c
    if ($isTextCol) {
        $cell2.NumberFormat = "@"
        $cell3.NumberFormat = "@"
    }

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2

    if ($isTextCol) {
        $cell2.ClearFormats()
        $cell3.ClearFormats()
    }
}

Write-Host "done"
